# Insert a new data row before row 50 (shifts existing rows 50-59 down to 51-60)
# and populate it with the new "Winter Nelis / Región del Maule" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 50 - this pushes rows 50..59 down to 51..60
# and automatically grows the sheet dimension to A1:T60.
$ws.Rows.Item(50).Insert()

# Fill in the values for the newly inserted row 50.
$ws.Cells.Item(50, 1).Value = 1
$ws.Cells.Item(50, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value = 45034
$ws.Cells.Item(50, 5).Value = 15
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100104
$ws.Cells.Item(50, 8).Value = "Frutos de pepita"
$ws.Cells.Item(50, 9).Value = 100104005
$ws.Cells.Item(50, 10).Value = "Pera"
$ws.Cells.Item(50, 11).Value = "Winter Nelis"
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 300
$ws.Cells.Item(50, 14).Value = 20000
$ws.Cells.Item(50, 15).Value = 21000
$ws.Cells.Item(50, 16).Value = 20500
$ws.Cells.Item(50, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(50, 18).Value = "Región del Maule"
$ws.Cells.Item(50, 19).Value = 1139
$ws.Cells.Item(50, 20).Value = 18
